$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add note about class being canceled on 2/17 due to Texas winter storm
$ws.Range("C13").Value = "**class canceled on 2/17 due to Texas winter storm"

# Push all remaining due dates back by one week (class was canceled 2/17/2021)
$ws.Range("B14").Value = 44251
$ws.Range("B16").Value = 44258
$ws.Range("B19").Value = 44265
$ws.Range("B21").Value = 44272
$ws.Range("B24").Value = 44279
$ws.Range("B26").Value = 44286
$ws.Range("B29").Value = 44293
$ws.Range("B31").Value = 44300
$ws.Range("B34").Value = 44307
$ws.Range("B36").Value = 44314

# Update the active selection to reflect where the editor left off
$ws.Range("B38").Select()
